# edit.ps1 - apply the "456a3b4" data refresh to 南宁-漫展信息.xlsx
#
# Summary of the change (derived from the OOXML diff):
#  - Sheet "展览" (展览, index 1): attendance/price counters (F, G columns)
#    refreshed for the 9 existing events; no rows added/removed.
#  - Sheet "演出" (演出, index 2): the "胡桃夹子" (Nutcracker) event (old
#    row 2) has dropped out of the feed. Every following row's B:I content
#    shifts up into the row above it, and the now-duplicate trailing row is
#    removed. Column A (the running index) is left untouched throughout.
#  - Sheet "本地生活" (index 3): untouched (still just the header row).
#  - Sheet "全部类型" (index 4): same "胡桃夹子" drop + up-shift as 演出,
#    plus the same attendance/price refresh as 展览 (applied at the rows'
#    *new*, post-shift positions).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition) - refresh F (want-to-go count) / G (price)
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F2").Value = 5493
$wsExpo.Range("F3").Value = 600
$wsExpo.Range("F4").Value = 12092
$wsExpo.Range("G4").Value = 62
$wsExpo.Range("F5").Value = 297
$wsExpo.Range("F6").Value = 610
$wsExpo.Range("F7").Value = 179
$wsExpo.Range("F8").Value = 325
$wsExpo.Range("F9").Value = 1102

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performance) - drop the first event (row 2), shifting
# every later row's content (columns B:I) up by one, then delete the
# trailing duplicate row. Column A is intentionally left alone.
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)
$lastRowShow = 7
for ($r = 2; $r -lt $lastRowShow; $r++) {
    $srcRow = $r + 1
    $src = $wsShow.Range("B" + $srcRow + ":I" + $srcRow)
    $dst = $wsShow.Range("B" + $r + ":I" + $r)
    $src.Copy($dst)
}
$wsShow.Rows.Item($lastRowShow).Delete()

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life) - no change.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - same drop/shift as 演出 ...
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$lastRowAll = 16
for ($r = 2; $r -lt $lastRowAll; $r++) {
    $srcRow = $r + 1
    $src = $wsAll.Range("B" + $srcRow + ":I" + $srcRow)
    $dst = $wsAll.Range("B" + $r + ":I" + $r)
    $src.Copy($dst)
}
$wsAll.Rows.Item($lastRowAll).Delete()

# ... plus the same attendance/price refresh as 展览, applied at the
# rows' new (post-shift) positions.
$wsAll.Range("F3").Value = 5493
$wsAll.Range("F4").Value = 600
$wsAll.Range("F6").Value = 12092
$wsAll.Range("G6").Value = 62
$wsAll.Range("F7").Value = 297
$wsAll.Range("F8").Value = 610
$wsAll.Range("F9").Value = 179
$wsAll.Range("F12").Value = 325
$wsAll.Range("F13").Value = 1102
